# summer 24 week 10 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.28
$ws.Range("B3").Value = 1.56
$ws.Range("F3").Value = 1.19
$ws.Range("C4").Value = 1.44
$ws.Range("F4").Value = 1.1
$ws.Range("C6").Value = 1.51
$ws.Range("D6").Value = 1.51
$ws.Range("F7").Value = 1.48
